# Update cryptos list data (prices and 1h volume/change percentages)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.820.98"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.635.16"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5083"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06428"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.43"
$ws.Range("E10").Value = "  +5.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07791"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.265"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.639.55"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "1.861.37"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5605"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "0.0₅7670"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "25.847.99"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.392"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.08"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.953"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.154"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.758"
$ws.Range("E25").Value = "  -6.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.31"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1234"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.850"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.239"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04967"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.309"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.252"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.568"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.386"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9029"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5567"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "1.133.81"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01569"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9962"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.17"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8008"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "0.0₈113"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.54"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4263"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.792"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05046"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9996"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.03%  "

Write-Host "Updated cryptos list"
